# Daily update at 8 AM UTC — append the next day's row of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 46

# Copy the formatting of the previous data row down onto the new row first,
# so the new cells inherit the same styles (e.g. the date format on column A).
$ws.Range("A45:D45").Copy()
$ws.Range("A46:D46").PasteSpecial(-4122) # xlPasteFormats

# Now fill in the actual values for the new day's data.
$ws.Cells.Item($newRow, 1).Value = 45995
$ws.Cells.Item($newRow, 2).Value = 104
$ws.Cells.Item($newRow, 3).Value = 115
$ws.Cells.Item($newRow, 4).Value = 109
